$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Classi" (sheet1): add row 13 for the week "28 marzo - 2 aprile"
# ---------------------------------------------------------------------------
$wsClassi = $wb.Worksheets.Item("Classi")

$wsClassi.Cells(13,1).Value = "28 marzo - 2 aprile"
$wsClassi.Cells(13,2).Value = 5393
$wsClassi.Cells(13,3).Value = 8157
$wsClassi.Cells(13,4).Value = 0.6509999999999999
$wsClassi.Cells(13,5).Value = 376539
$wsClassi.Cells(13,6).Value = 244617
$wsClassi.Cells(13,7).Value = 0.65
$wsClassi.Cells(13,8).Value = 244530
$wsClassi.Cells(13,9).Value = 18852
$wsClassi.Cells(13,10).Value = 1
$wsClassi.Cells(13,11).Value = 0.077
$wsClassi.Cells(13,14).Value = 87
$wsClassi.Cells(13,15).Value = 0.001

$wsClassi.Range("B13:C13").NumberFormat = "#,##0"
$wsClassi.Range("E13:F13").NumberFormat = "#,##0"
$wsClassi.Range("H13:I13").NumberFormat = "#,##0"
$wsClassi.Range("D13").NumberFormat = "0.0%"
$wsClassi.Range("G13").NumberFormat = "0.0%"
$wsClassi.Range("J13:K13").NumberFormat = "0.0%"

# O13 reuses the exact same style as O8:O12 (0.0% with the alternate font) -
# copy the formatting only from O12 so the existing style entry is reused.
$wsClassi.Cells(12,15).Copy() | Out-Null
$wsClassi.Cells(13,15).PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------------
# Sheet "Alunni in presenza" (sheet2): add row 13
# ---------------------------------------------------------------------------
$wsPresenza = $wb.Worksheets.Item("Alunni in presenza")

$wsPresenza.Cells(13,1).Value = "28 marzo - 2 aprile"
$wsPresenza.Cells(13,2).Value = 7395000
$wsPresenza.Cells(13,3).Value = 4792852
$wsPresenza.Cells(13,4).Value = 0.648
$wsPresenza.Cells(13,5).Value = 4655153
$wsPresenza.Cells(13,6).Value = 0.971

$wsPresenza.Range("B13:C13").NumberFormat = "#,##0"
$wsPresenza.Range("E13").NumberFormat = "#,##0"
$wsPresenza.Range("D13").NumberFormat = "0.0%"
$wsPresenza.Range("F13").NumberFormat = "0.0%"

# ---------------------------------------------------------------------------
# Sheet "Alunni" (sheet3): add rows 35, 36, 37 (Infanzia / Primaria / Sec.)
# ---------------------------------------------------------------------------
$wsAlunni = $wb.Worksheets.Item("Alunni")

$wsAlunni.Cells(35,1).Value = "28 marzo - 2 aprile"
$wsAlunni.Cells(35,2).Value = "Infanzia"
$wsAlunni.Cells(35,3).Value = 547613
$wsAlunni.Cells(35,4).Value = 536432
$wsAlunni.Cells(35,5).Value = 11181
$wsAlunni.Cells(35,6).Value = 0.02

$wsAlunni.Cells(36,1).Value = "28 marzo - 2 aprile"
$wsAlunni.Cells(36,2).Value = "Primaria"
$wsAlunni.Cells(36,3).Value = 1513526
$wsAlunni.Cells(36,4).Value = 1464519
$wsAlunni.Cells(36,5).Value = 49007
$wsAlunni.Cells(36,6).Value = 0.032

$wsAlunni.Cells(37,1).Value = "28 marzo - 2 aprile"
$wsAlunni.Cells(37,2).Value = "Sec. 1° e 2° Grado"
$wsAlunni.Cells(37,3).Value = 2731713
$wsAlunni.Cells(37,4).Value = 2654202
$wsAlunni.Cells(37,5).Value = 77511
$wsAlunni.Cells(37,6).Value = 0.027999999999999997

$wsAlunni.Range("C35:E37").NumberFormat = "#,##0"
$wsAlunni.Range("F35:F37").NumberFormat = "0.0%"

# ---------------------------------------------------------------------------
# Sheet "Personale scolastico" (sheet4): add row 13
# ---------------------------------------------------------------------------
$wsPersonale = $wb.Worksheets.Item("Personale scolastico")

$wsPersonale.Cells(13,1).Value = "28 marzo - 2 aprile"
$wsPersonale.Cells(13,2).Value = 775867
$wsPersonale.Cells(13,3).Value = 500681
$wsPersonale.Cells(13,4).Value = 0.645
$wsPersonale.Cells(13,5).Value = 475864
$wsPersonale.Cells(13,6).Value = 0.95
$wsPersonale.Cells(13,7).Value = 204526
$wsPersonale.Cells(13,8).Value = 132766
$wsPersonale.Cells(13,9).Value = 0.649
$wsPersonale.Cells(13,10).Value = 127749
$wsPersonale.Cells(13,11).Value = 0.9620000000000001

$wsPersonale.Range("B13:C13").NumberFormat = "#,##0"
$wsPersonale.Range("E13").NumberFormat = "#,##0"
$wsPersonale.Range("G13:H13").NumberFormat = "#,##0"
$wsPersonale.Range("J13").NumberFormat = "#,##0"
$wsPersonale.Range("D13").NumberFormat = "0.0%"
$wsPersonale.Range("F13").NumberFormat = "0.0%"
$wsPersonale.Range("I13").NumberFormat = "0.0%"
$wsPersonale.Range("K13").NumberFormat = "0.0%"

# ---------------------------------------------------------------------------
# Restore the selection on every sheet to match the post-edit state, then
# re-select on the sheet that must remain the active tab (Personale
# scolastico) so it keeps focus.
# ---------------------------------------------------------------------------
$wsClassi.Range("O14").Select() | Out-Null
$wsPresenza.Range("B14").Select() | Out-Null
$wsAlunni.Range("C38").Select() | Out-Null
$wsPersonale.Range("B14").Select() | Out-Null
